# Populate instruction, time, numberOfIngredients, and ranking columns
# for recipes 3-10 (rows 4-11), matching the wrapText style already used
# by rows 2-3 (cellXfs index 1: default font + wrapText).
#
# NOTE: row 5 is written before row 4 so the new shared-string table
# entries land in the same order as the target workbook (index 18 =
# the Stir-fried Egg and Tomato instruction referenced by row 5, index
# 19 = the Sun-dried Tomato Omelet instruction referenced by row 4).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$s18 = @'
Beat eggs with 1/2 teaspoon salt until smooth but not frothy.
Heat 1 tablespoon oil in a 12-inch nonstick skillet over medium-high heat until hot. Add eggs and cook, undisturbed, just until a thin film of cooked egg forms on bottom of skillet but most of eggs are still runny, 5 to 10 seconds. Immediately scrape eggs into a bowl. Wipe out skillet.
Heat remaining tablespoon oil in skillet over medium-high heat until hot. Add scallions and stir-fry until just softened, about 30 seconds. Add tomatoes and cook, stirring and turning occasionally, until juices are released and tomatoes are slightly wilted but still intact, 4 to 6 minutes. Sprinkle sugar and 1/4 teaspoon salt over tomatoes and stir to combine. Return eggs to skillet and cook, stirring occasionally, until eggs are just cooked through. Serve sprinkled with reserved scallion greens.
'@

$s19 = @'
In a medium bowl, beat the eggs with a pinch each of salt and pepper.
'@

$s20 = @'
Combine the carrots, potatoes, onion, cabbage, garlic, chicken stock, olive oil, thyme, basil, parsley, salt, and pepper in a stock pot over medium-high heat; bring to a simmer and cook until the carrots are tender, about 20 minutes. Transfer to a blender in small batches and blend until smooth.
'@

$s21 = @'
For the dressing: Place the anchovies into a blender or food processor. Throw in the Dijon mustard, vinegar, Worcestershire, garlic and lemon juice. Pulse the processor or blend on low speed for several seconds. Scrape down the sides.
With the food processor or blender on, drizzle the olive oil into the mixture in a small stream. Scrape down the sides. Add the Parmesan, salt and a generous grind of black pepper. Pulse the whole thing together and mix until thoroughly combined. Refrigerate the dressing for a few hours (it just gets better!) before using it on the salad.
For the croutons: Slice the bread into thick slices and cut them into 1-inch cubes. Throw them onto a baking sheet.
Heat the olive oil in a small saucepan or skillet over low heat.
Crush-but don't chop-the garlic and add them to the oil. Use a spoon to move the garlic around in the pan. After 3 to 5 minutes, turn off the heat and remove the garlic from the pan.
Slowly drizzle the olive oil over the bread cubes. Mix together with your hands, and then sprinkle lightly with salt. Toss and cook in the pan until golden brown and crisp. (Add a little butter for more flavor!)
For the salad: Wash and dry the hearts of romaine lettuce. Leave them whole. Use a vegetable peeler and shave off large, thin slices of Parmesan.
Drizzle about half of the dressing over the top of the hearts. Throw in a good handful of the Parmesan shavings. Give it a good initial toss, just so you can evaluate how much more dressing you need.
Add more dressing and Parmesan to taste. Add the cooled croutons. Toss gently.
Read more at: http://www.foodnetwork.com/recipes/ree-drummond/caesar-salad-recipe.html?oc=linkback
'@

$s22 = @'
Preheat the oven to 350 degrees F. Spray a baking sheet with nonstick cooking spray.
Mix together the flour and 1 teaspoon of the House Seasoning in a small bowl. Sprinkle the chicken with the remaining 1 teaspoon House Seasoning. Pour the buttermilk into a shallow dish. Dredge the chicken in the buttermilk, followed by the flour.
Heat 2 tablespoons of the oil in a heavy nonstick skillet over medium-high heat. Add half of the chicken breasts to the hot oil and cook until both sides are browned, about 3 minutes per side, and then transfer to the baking sheet. Repeat with the remaining chicken breasts. Transfer the baking sheet to the oven and bake until the chicken is cooked through, about 10 minutes.
Whisk together the cornstarch and 1/4 cup chicken broth until dissolved. Set aside.
To make the gravy , add the remaining 1 tablespoon oil to the same skillet and heat over medium heat. Add the onions and saute until translucent, about 2 minutes.
Add the remaining chicken broth, scrape the pan drippings with a wooden spoon, raise the heat to medium-high and cook until the mixture begins to bubble, about 2 minutes. Stir in the dissolved cornstarch to incorporate. Bring to a simmer and continue to cook until the mixture thickens, 4 to 5 minutes. Stir in the milk and black pepper and continue cooking over medium-high heat until thickened, 5 minutes longer.
Remove the chicken from the oven and top each piece with 4 teaspoons of the gravy. Sprinkle with the chopped green onions.
Mix together the salt, garlic powder and pepper.
'@

$s23 = @'
Line up 4 parfait, white wine, or other tall glasses.
Spoon 2 tablespoons of yogurt into each glass and smooth surface.
Spoon 2 tablespoons of granola overtop and smooth surface.
Spoon 2 tablespoons of fruit overtop and smooth surface.
Repeat the process, adding a bit of honey here and there, to taste.
Read more at: http://www.foodnetwork.com/recipes/granola-yogurt-berry-parfait-recipe.html?oc=linkback
'@

$s24 = @'
Cook the bacon in a large, deep skillet over medium-high heat until evenly browned, about 10 minutes. Drain the bacon slices on a paper towel-lined plate.
Arrange the cooked bacon, lettuce, and tomato slices on one slice of bread. Spread one side of remaining bread slice with the mayonnaise. Bring the two pieces together to make a sandwich.
'@

$s25 = @'
Preheat skillet over medium heat. Generously butter one side of a slice of bread. Place bread butter-side-down onto skillet bottom and add 1 slice of cheese. Butter a second slice of bread on one side and place butter-side-up on top of sandwich. Grill until lightly browned and flip over; continue grilling until cheese is melted. Repeat with remaining 2 slices of bread, butter and slice of cheese.
'@

$ws.Rows(5).RowHeight = 409.5
$ws.Range("C5").WrapText = $true
$ws.Range("C5").Value = $s18
$ws.Range("D5").Value = 15
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 4

$ws.Rows(4).RowHeight = 120
$ws.Range("C4").WrapText = $true
$ws.Range("C4").Value = $s19
$ws.Range("D4").Value = 25
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = 5

$ws.Rows(6).RowHeight = 409.5
$ws.Range("C6").WrapText = $true
$ws.Range("C6").Value = $s20
$ws.Range("D6").Value = 50
$ws.Range("E6").Value = 10
$ws.Range("F6").Value = 5

$ws.Rows(7).RowHeight = 409.5
$ws.Range("C7").WrapText = $true
$ws.Range("C7").Value = $s21
$ws.Range("D7").Value = 20
$ws.Range("E7").Value = 10
$ws.Range("F7").Value = 5

$ws.Rows(8).RowHeight = 409.5
$ws.Range("C8").WrapText = $true
$ws.Range("C8").Value = $s22
$ws.Range("D8").Value = 45
$ws.Range("E8").Value = 10
$ws.Range("F8").Value = 5

$ws.Rows(9).RowHeight = 409.5
$ws.Range("C9").WrapText = $true
$ws.Range("C9").Value = $s23
$ws.Range("D9").Value = 5
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = 5

$ws.Rows(10).RowHeight = 409.5
$ws.Range("C10").WrapText = $true
$ws.Range("C10").Value = $s24
$ws.Range("D10").Value = 15
$ws.Range("E10").Value = 5
$ws.Range("F10").Value = 5

$ws.Rows(11).RowHeight = 409.5
$ws.Range("C11").WrapText = $true
$ws.Range("C11").Value = $s25
$ws.Range("D11").Value = 20
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 5

# Reset the sheet view: clear the frozen/scrolled top-left cell and move
# the active selection to B1.
$ws.Range("B1").Select()

